# Update calibrated cost values for selected rows.
# Columns J (10) through AS (45) hold the same constant value across the row
# for each of the affected variables; only the constant changes per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("strategy_id-0")

$updates = @{
    100 = 401896.5761
    101 = 73917.89659999999
    102 = 872582.0111
    103 = 13256.58684
    104 = 309151.2124
    105 = 30141.22701
    106 = 95984.02122
    107 = 215099.4036
    114 = 50.40808688
    115 = 14305407.6
}

$firstCol = 10  # J
$lastCol  = 45  # AS

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $rng = $ws.Range($ws.Cells.Item($row, $firstCol), $ws.Cells.Item($row, $lastCol))
    $rng.Value = $value
}
